$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (A1) now reuses the Kyrgyz 4.2.2 title string ----------
$ws.Range("A1").Value = "4.2.2 Уюштурулган окутуулардын түрлөрүнө катышуунун деңгээли (мектепке кире турган расмий жаш куракка чыкканга чейинки бир жыл үчүн)"

# --- New column R (2023) -------------------------------------------------
$ws.Cells.Item(4, 18).Value = 2023
$ws.Cells.Item(5, 18).Value = 53.5

# --- Row 5 now carries country names instead of duplicated titles --------
$ws.Range("A5").Value = "Кыргыз Республикасы"
$ws.Range("B5").Value = "Кыргызская Республика"
$ws.Range("C5").Value = "Kyrgyz Republic "

# --- Row height on row 5 shrinks now that it's a short label -------------
$ws.Rows("5:5").RowHeight = 21

# --- Column widths: columns A:C share one uniform width -------------------
$ws.Columns("A:C").ColumnWidth = 35.85546875
